$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4269
$ws.Range("I62").Value = 1690
$ws.Range("J62").Value = 4784.8
$ws.Range("K62").Value = 1690
$ws.Range("L62").Value = 4784.8
$ws.Range("M62").Value = -1066
$ws.Range("N62").Value = -6032.8

$ws.Range("H65").Value = 4269
$ws.Range("I65").Value = 1690
$ws.Range("J65").Value = 4784.8
$ws.Range("K65").Value = 8450
$ws.Range("L65").Value = 23924
$ws.Range("M65").Value = -5330
$ws.Range("N65").Value = -30164

$ws.Range("H112").Value = 45456264
$ws.Range("I112").Value = 645
$ws.Range("J112").Value = 55557510
$ws.Range("K112").Value = 1935
$ws.Range("L112").Value = 166672530
$ws.Range("M112").Value = -827
$ws.Range("N112").Value = -166674746

$ws.Range("H132").Value = 235164.75
$ws.Range("I132").Value = 240735.33
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 722205.99
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -719675.99
$ws.Range("N132").Value = -8660

$ws.Range("H141").Value = 1354.591
$ws.Range("I141").Value = 516.73334
$ws.Range("J141").Value = 3150
$ws.Range("K141").Value = 1550.20002
$ws.Range("L141").Value = 9450
$ws.Range("M141").Value = 3629.79998
$ws.Range("N141").Value = -19810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 400
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -384

$ws.Range("H5").Value = 183.5
$ws.Range("I5").Value = 80
$ws.Range("J5").Value = 204.2
$ws.Range("K5").Value = 80
$ws.Range("L5").Value = 204.2
$ws.Range("M5").Value = 32
$ws.Range("N5").Value = -428.2

$ws.Range("H45").Value = 1621.3334
$ws.Range("I45").Value = 1723
$ws.Range("J45").Value = 1540
$ws.Range("K45").Value = 1723
$ws.Range("L45").Value = 1540
$ws.Range("M45").Value = -1346
$ws.Range("N45").Value = -2294

$ws.Range("H74").Value = 4033.4878
$ws.Range("I74").Value = 942.9091
$ws.Range("J74").Value = 7612.0527
$ws.Range("K74").Value = 942.9091
$ws.Range("L74").Value = 7612.0527
$ws.Range("M74").Value = -68.90909999999997
$ws.Range("N74").Value = -9360.0527

$ws.Range("H77").Value = 4033.4878
$ws.Range("I77").Value = 942.9091
$ws.Range("J77").Value = 7612.0527
$ws.Range("K77").Value = 4714.5455
$ws.Range("L77").Value = 38060.2635
$ws.Range("M77").Value = -346.5455000000002
$ws.Range("N77").Value = -46796.2635

$ws.Range("H102").Value = 2066.7
$ws.Range("I102").Value = 2029.8572
$ws.Range("J102").Value = 2152.6667
$ws.Range("K102").Value = 2029.8572
$ws.Range("L102").Value = 2152.6667
$ws.Range("M102").Value = -407.8571999999999
$ws.Range("N102").Value = -5396.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 183.5
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 204.2
$ws.Range("K4").Value = 80
$ws.Range("L4").Value = 204.2
$ws.Range("M4").Value = 35
$ws.Range("N4").Value = -434.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2032121
$ws.Range("I6").Value = 2540000.8
$ws.Range("J6").Value = 602
$ws.Range("K6").Value = 2540000.8
$ws.Range("L6").Value = 602
$ws.Range("M6").Value = -2539887.8
$ws.Range("N6").Value = -828

$ws.Range("H7").Value = 100.25
$ws.Range("I7").Value = 100.333336
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 100.333336
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 12.666664
$ws.Range("N7").Value = -326

$ws.Range("H132").Value = 2234
$ws.Range("I132").Value = 2194.8333
$ws.Range("J132").Value = 2351.5
$ws.Range("K132").Value = 6584.499899999999
$ws.Range("L132").Value = 7054.5
$ws.Range("M132").Value = -4054.499899999999
$ws.Range("N132").Value = -12114.5

$ws.Range("H134").Value = 2656.8667
$ws.Range("I134").Value = 3569.3333
$ws.Range("J134").Value = 1288.1666
$ws.Range("K134").Value = 10707.9999
$ws.Range("L134").Value = 3864.4998
$ws.Range("M134").Value = -8172.999899999999
$ws.Range("N134").Value = -8934.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6585657
$ws.Range("I80").Value = 14486057
$ws.Range("J80").Value = 1990
$ws.Range("K80").Value = 43458171
$ws.Range("L80").Value = 5970
$ws.Range("M80").Value = -43457235
$ws.Range("N80").Value = -7842

$ws.Range("H83").Value = 6585657
$ws.Range("I83").Value = 14486057
$ws.Range("J83").Value = 1990
$ws.Range("K83").Value = 130374513
$ws.Range("L83").Value = 17910
$ws.Range("M83").Value = -130369833
$ws.Range("N83").Value = -27270

$ws.Range("H129").Value = 903.5625
$ws.Range("I129").Value = 614.1111
$ws.Range("K129").Value = 1842.3333
$ws.Range("M129").Value = 3157.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 620
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 833.3333
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 833.3333
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -1423.3333

$ws.Range("H27").Value = 620
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 833.3333
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 833.3333
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -1047.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 283.16666
$ws.Range("I107").Value = 182.5
$ws.Range("J107").Value = 383.83334
$ws.Range("K107").Value = 547.5
$ws.Range("L107").Value = 1151.50002
$ws.Range("M107").Value = 1372.5
$ws.Range("N107").Value = -4991.500019999999

$ws.Range("H109").Value = 47938
$ws.Range("J109").Value = 47938
$ws.Range("L109").Value = 47938
$ws.Range("N109").Value = -50712

$ws.Range("H136").Value = 10186.471
$ws.Range("I136").Value = 12395.407
$ws.Range("J136").Value = 1666.2858
$ws.Range("K136").Value = 37186.221
$ws.Range("L136").Value = 4998.857400000001
$ws.Range("M136").Value = -34636.221
$ws.Range("N136").Value = -10098.8574
